# New PO forecast model
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Weekly Quantity": append a new week row (row 4)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A4").Value = 45676.99999999999
$wsWeekly.Range("A4").NumberFormat = $wsWeekly.Range("A3").NumberFormat
$wsWeekly.Range("B4").Value = 8

# ---------------------------------------------------------------------
# Sheet "Monthly Trend": append a new month row (row 4)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A4").Value = 45688.99999999999
$wsMonthly.Range("A4").NumberFormat = $wsMonthly.Range("A3").NumberFormat
$wsMonthly.Range("B4").Value = 8

# ---------------------------------------------------------------------
# Sheet "PO Forecast": refresh the forecast values/dates and extend one
# more row into the future
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$wsForecast.Range("B2").Value = 1500
$wsForecast.Range("B3").Value = 449

$wsForecast.Range("A4").Value = 45676.99999999999
$wsForecast.Range("A5").Value = 45683.99999999999
$wsForecast.Range("A6").Value = 45690.99999999999
$wsForecast.Range("A7").Value = 45697.99999999999
$wsForecast.Range("A8").Value = 45704.99999999999
$wsForecast.Range("A9").Value = 45711.99999999999
$wsForecast.Range("A10").Value = 45718.99999999999
$wsForecast.Range("A11").Value = 45725.99999999999

$wsForecast.Range("A12").Value = 45732.99999999999
$wsForecast.Range("A12").NumberFormat = $wsForecast.Range("A11").NumberFormat
$wsForecast.Range("B12").Value = 0
